# Apply "Natmi following Dr Hou advice" edit: recompute the LR-pairs stats
# including an "ECs" sending-cluster block, over sheet1 rows 2-17 (A1:T17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("ECs", "Serpinf1", "Plxdc2", "ECs", 2, 0.6666666666666666, 5.018449, 15.055347, 0.007216928732731275, 0.007216928732731275, 3, 1, 4.416526333333333, 13.249579, 0.03380457256708, 0.03380457256708, 22.16411216099033, 199.477009448913, 0.0002439651910570591, 0.0002439651910570591)
    ,("ECs", "Serpinf1", "Plxdc2", "FAPs", 2, 0.6666666666666666, 5.018449, 15.055347, 0.007216928732731275, 0.007216928732731275, 3, 1, 111.74676, 335.24028, 0.8553218462766415, 0.8553218462766417, 560.79541597524, 5047.15874377716, 0.006172796808126657, 0.006172796808126658)
    ,("ECs", "Serpinf1", "Plxdc2", "M2", 2, 0.6666666666666666, 5.018449, 15.055347, 0.007216928732731275, 0.007216928732731275, 3, 1, 3.275441, 9.826323, 0.02507058140647845, 0.02507058140647845, 16.437633611009, 147.938702499081, 0.0001809325992986928, 0.0001809325992986928)
    ,("ECs", "Serpinf1", "Plxdc2", "sCs", 2, 0.6666666666666666, 5.018449, 15.055347, 0.007216928732731275, 0.007216928732731275, 3, 1, 11.21005766666667, 33.630173, 0.08580299974979994, 0.08580299974979995, 56.25710268722567, 506.313924185031, 0.0006192341342488656, 0.0006192341342488657)
    ,("FAPs", "Serpinf1", "Plxdc2", "ECs", 3, 1, 675.4945473333333, 2026.483642, 0.9714148748852961, 0.9714148748852962, 3, 1, 4.416526333333333, 13.249579, 0.03380457256708, 0.03380457256708, 2983.339456320746, 26850.05510688672, 0.03283826463080093, 0.03283826463080093)
    ,("FAPs", "Serpinf1", "Plxdc2", "FAPs", 3, 1, 675.4945473333333, 2026.483642, 0.9714148748852961, 0.9714148748852962, 3, 1, 111.74676, 335.24028, 0.8553218462766415, 0.8553218462766417, 75484.32706216663, 679358.9435594998, 0.8308723642874842, 0.8308723642874845)
    ,("FAPs", "Serpinf1", "Plxdc2", "M2", 3, 1, 675.4945473333333, 2026.483642, 0.9714148748852961, 0.9714148748852962, 3, 1, 3.275441, 9.826323, 0.02507058140647845, 0.02507058140647845, 2212.542535612041, 19912.88282050837, 0.02435393570027589, 0.02435393570027589)
    ,("FAPs", "Serpinf1", "Plxdc2", "sCs", 3, 1, 675.4945473333333, 2026.483642, 0.9714148748852961, 0.9714148748852962, 3, 1, 11.21005766666667, 33.630173, 0.08580299974979994, 0.08580299974979995, 7572.332829125562, 68150.99546213007, 0.08335031026673499, 0.08335031026673502)
    ,("M2", "Serpinf1", "Plxdc2", "ECs", 3, 1, 1.832754333333333, 5.498263, 0.002635646473297045, 0.002635646473297045, 3, 1, 4.416526333333333, 13.249579, 0.03380457256708, 0.03380457256708, 8.094407775697443, 72.84966998127699, 0.00008909690246773845, 0.00008909690246773845)
    ,("M2", "Serpinf1", "Plxdc2", "FAPs", 3, 1, 1.832754333333333, 5.498263, 0.002635646473297045, 0.002635646473297045, 3, 1, 111.74676, 335.24028, 0.8553218462766415, 0.8553218462766417, 204.80435862596, 1843.23922763364, 0.002254326007672948, 0.002254326007672948)
    ,("M2", "Serpinf1", "Plxdc2", "M2", 3, 1, 1.832754333333333, 5.498263, 0.002635646473297045, 0.002635646473297045, 3, 1, 3.275441, 9.826323, 0.02507058140647845, 0.02507058140647845, 6.003078686327666, 54.027708176949, 0.00006607718946749139, 0.00006607718946749139)
    ,("M2", "Serpinf1", "Plxdc2", "sCs", 3, 1, 1.832754333333333, 5.498263, 0.002635646473297045, 0.002635646473297045, 3, 1, 11.21005766666667, 33.630173, 0.08580299974979994, 0.08580299974979995, 20.54528176549989, 184.907535889499, 0.0002261463736888675, 0.0002261463736888675)
    ,("sCs", "Serpinf1", "Plxdc2", "ECs", 3, 1, 13.02608766666667, 39.078263, 0.0187325499086756, 0.0187325499086756, 3, 1, 4.416526333333333, 13.249579, 0.03380457256708, 0.03380457256708, 57.53005920014188, 517.770532801277, 0.000633245842754272, 0.000633245842754272)
    ,("sCs", "Serpinf1", "Plxdc2", "FAPs", 3, 1, 13.02608766666667, 39.078263, 0.0187325499086756, 0.0187325499086756, 3, 1, 111.74676, 335.24028, 0.8553218462766415, 0.8553218462766417, 1455.62309222596, 13100.60783003364, 0.01602235917335774, 0.01602235917335775)
    ,("sCs", "Serpinf1", "Plxdc2", "M2", 3, 1, 13.02608766666667, 39.078263, 0.0187325499086756, 0.0187325499086756, 3, 1, 3.275441, 9.826323, 0.02507058140647845, 0.02507058140647845, 42.66618161299434, 383.995634516949, 0.000469635917436372, 0.000469635917436372)
    ,("sCs", "Serpinf1", "Plxdc2", "sCs", 3, 1, 13.02608766666667, 39.078263, 0.0187325499086756, 0.0187325499086756, 3, 1, 11.21005766666667, 33.630173, 0.08580299974979994, 0.08580299974979995, 146.0231939143888, 1314.208745229499, 0.001607308975127207, 0.001607308975127207)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}